# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it is used
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2) Narrow the now-shorter "Status" columns:
#    Overview columns E & F, and column C on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status values -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus

# --- Narrow the status columns to fit the shorter text -----------------------
# Target stored column width is 13.4101848602295 "characters" worth of pixels.
# The engine (like real Excel) snaps stored widths to whole pixels, so the
# ColumnWidth we *assign* is chosen to land on the closest achievable stored
# value.
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede.Columns.Item(3).ColumnWidth = $newWidth
